$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matches the source data, where numeric-looking IDs/years/CCCD numbers are
# kept as text so leading zeros etc. survive). We briefly flip the cell to a
# text number-format, assign the value, then clear the format again so the
# workbook's style table ends up the same as before the edit.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2: "huy" -> "Huy", CCCD number replaced, Số Người fixed to 3
Set-TextValue 2 2 "Huy"
Set-TextValue 2 4 "042203013460"
Set-TextValue 2 5 "3"

# Row 3: becomes an exact copy of the (new) row 2 data
Set-TextValue 3 1 "1"
Set-TextValue 3 2 "Huy"
Set-TextValue 3 3 "2003"
Set-TextValue 3 4 "042203013460"
Set-TextValue 3 5 "3"

# Row 4 (new)
Set-TextValue 4 1 "3"
Set-TextValue 4 2 "Huy"
Set-TextValue 4 3 "2003"
Set-TextValue 4 4 "042203013460"
Set-TextValue 4 5 "3"

# Row 5 (new)
Set-TextValue 5 1 "3"
Set-TextValue 5 2 "Huy"
Set-TextValue 5 3 "2003"
Set-TextValue 5 4 "042203013460"
Set-TextValue 5 5 "3"
